$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 holds the "C25" capacitor (560p). It is being changed to the same
# 100pF part already used by "C1" (row 2), so the two BoM lines are merged:
# C1's quantity becomes 2 and its reference list becomes "C1 C25", while the
# now-duplicate C25 row is removed entirely (rows below shift up by one).

$ws.Range("A2").Value = "C1 C25"
$ws.Range("B2").Value = 2

$ws.Rows.Item(9).Delete()

$excel.CalculateFullRebuild()
